$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$src = $ws.Range("A98:H98")
$dst = $ws.Range("A98:H101")
$src.AutoFill($dst)
Write-Host "UsedRange:" $ws.UsedRange.Address
